$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A276").Formula = "=+A275+1"
$ws.Range("B276").Value = "VIỄN THÔNG NINH BÌNH  TẬP ĐOÀN BƯU CHÍNH VIỄN THÔNG VIỆT NAM"
$ws.Range("C276").Value = "2700141673"
$ws.Range("D276").Value = "KHTV.170"

$ws.Range("A277").Formula = "=+A276+1"
$ws.Range("B277").Value = "TỔNG CÔNG TY CỔ PHẦN BƯU CHÍNH VIETTEL"
$ws.Range("C277").Value = "0104093672"
$ws.Range("D277").Value = "KHTV.176"

$ws.Range("A278").Formula = "=+A277+1"
$ws.Range("D278").Value = "KHTV.175"
$ws.Range("C278").Value = "0108461831"
$ws.Range("B278").Value = "CÔNG TY TNHH DỊCH VỤ DU LỊCH THỜI ĐẠI"

$ws.Range("A279").Formula = "=+A278+1"
$ws.Range("D279").Value = "KHTV.172"
$ws.Range("B279").Value = "CÔNG TY CỔ PHẦN NAZ TECCON"
$ws.Range("C279").Value = "0601304965"
